$d = $word.ActiveDocument

$pairs = @(
    ,@("68×16=", "39×37=")
    ,@("77×10=", "57×73=")
    ,@("79×13=", "65×77=")
    ,@("84×73=", "49×23=")
    ,@("67×12=", "91×54=")
    ,@("12×35=", "55×94=")
    ,@("93×98=", "20×89=")
    ,@("79×69=", "89×39=")
    ,@("19×62=", "95×91=")
    ,@("21×68=", "99×71=")
    ,@("33×67=", "55×59=")
    ,@("85×19=", "62×59=")
    ,@("49×68=", "85×100=")
    ,@("81×28=", "14×71=")
    ,@("74×63=", "31×97=")
    ,@("59×49=", "23×83=")
    ,@("32×30=", "63×55=")
    ,@("35×67=", "51×94=")
    ,@("11×92=", "58×39=")
    ,@("95×74=", "87×43=")
    ,@("25×61=", "12×14=")
    ,@("12×16=", "29×82=")
    ,@("16×61=", "18×34=")
    ,@("57×44=", "25×44=")
    ,@("36×41=", "34×60=")
    ,@("90×77=", "98×94=")
    ,@("52×48=", "67×64=")
    ,@("10×57=", "59×77=")
    ,@("80×10=", "13×73=")
    ,@("15×90=", "27×46=")
    ,@("63×99=", "77×46=")
    ,@("51×51=", "24×13=")
    ,@("46×92=", "16×77=")
    ,@("12×21=", "68×73=")
    ,@("85×85=", "83×40=")
    ,@("14×53=", "34×74=")
    ,@("21×64=", "62×60=")
    ,@("57×51=", "79×58=")
    ,@("28×44=", "47×37=")
    ,@("10×90=", "56×45=")
    ,@("25×84=", "18×23=")
    ,@("56×22=", "71×42=")
    ,@("35×24=", "62×53=")
    ,@("100×46=", "31×27=")
    ,@("21×56=", "53×64=")
    ,@("63×97=", "10×89=")
    ,@("13×81=", "88×69=")
    ,@("43×53=", "69×88=")
    ,@("32×87=", "68×60=")
    ,@("46×36=", "35×69=")
    ,@("54×50=", "30×48=")
    ,@("38×76=", "82×28=")
    ,@("100×92=", "24×40=")
    ,@("30×71=", "13×90=")
    ,@("45×55=", "52×88=")
    ,@("49×36=", "79×25=")
    ,@("76×91=", "63×36=")
    ,@("76×66=", "89×69=")
    ,@("11×88=", "42×28=")
    ,@("20×73=", "63×81=")
    ,@("54×19=", "54×75=")
    ,@("82×60=", "88×15=")
    ,@("72×28=", "71×61=")
    ,@("26×24=", "90×67=")
    ,@("58×49=", "44×91=")
    ,@("90×12=", "94×44=")
    ,@("51×39=", "19×81=")
    ,@("77×74=", "96×13=")
    ,@("28×79=", "62×21=")
    ,@("69×27=", "32×84=")
    ,@("56×24=", "47×68=")
    ,@("11×32=", "100×97=")
    ,@("11×50=", "64×21=")
    ,@("95×72=", "23×39=")
    ,@("72×87=", "85×64=")
    ,@("28×80=", "16×29=")
    ,@("58×47=", "17×50=")
    ,@("93×76=", "44×78=")
    ,@("76×99=", "90×40=")
    ,@("29×46=", "33×88=")
    ,@("10×73=", "96×86=")
    ,@("62×76=", "54×18=")
    ,@("50×21=", "36×96=")
    ,@("71×15=", "58×21=")
    ,@("75×78=", "92×54=")
    ,@("90×47=", "42×63=")
    ,@("84×28=", "79×10=")
    ,@("86×38=", "76×31=")
    ,@("50×20=", "16×85=")
    ,@("63×96=", "74×87=")
    ,@("65×34=", "61×92=")
    ,@("32×59=", "27×57=")
    ,@("92×21=", "25×77=")
    ,@("32×85=", "72×69=")
    ,@("89×72=", "100×68=")
    ,@("73×75=", "66×42=")
    ,@("34×67=", "22×64=")
    ,@("42×94=", "67×25=")
    ,@("26×62=", "84×36=")
    ,@("67×75=", "97×38=")
)

foreach ($pair in $pairs) {
    $old = $pair[0]
    $new = $pair[1]
    $range = $d.Content
    $found = $range.Find.Execute($old, $true, $true, $false, $false, $false, $true, 1, $false, $new, 2)
    if (-not $found) {
        Write-Output "NOT FOUND: $old"
    }
}

Write-Output "Done"
